$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")
$ws.Range("A20").Value = 100
